$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 342
$ws.Range("I4").Value = 342
$ws.Range("K4").Value = 342
$ws.Range("M4").Value = -228
$ws.Range("H5").Value = 211.85715
$ws.Range("I5").Value = 233.5
$ws.Range("K5").Value = 233.5
$ws.Range("M5").Value = -118.5
$ws.Range("H6").Value = 2339.3845
$ws.Range("I6").Value = 41.2
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 123.6
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = -11.60000000000001
$ws.Range("N6").Value = -30224
$ws.Range("H9").Value = 20834928
$ws.Range("I9").Value = 25001896
$ws.Range("K9").Value = 25001896
$ws.Range("M9").Value = -25001727
$ws.Range("H17").Value = 345406.22
$ws.Range("J17").Value = 345406.22
$ws.Range("L17").Value = 1036218.66
$ws.Range("N17").Value = -1036554.66
$ws.Range("H18").Value = 2072.6667
$ws.Range("I18").Value = 2072.6667
$ws.Range("K18").Value = 2072.6667
$ws.Range("M18").Value = -1788.6667
$ws.Range("H53").Value = 92104.27
$ws.Range("I53").Value = 268.85715
$ws.Range("J53").Value = 252816.25
$ws.Range("K53").Value = 268.85715
$ws.Range("L53").Value = 252816.25
$ws.Range("M53").Value = 368.14285
$ws.Range("N53").Value = -254090.25
$ws.Range("H62").Value = 901515.9399999999
$ws.Range("I62").Value = 1472914
$ws.Range("J62").Value = 101558.6
$ws.Range("K62").Value = 1472914
$ws.Range("L62").Value = 101558.6
$ws.Range("M62").Value = -1472290
$ws.Range("N62").Value = -102806.6
$ws.Range("H65").Value = 901515.9399999999
$ws.Range("I65").Value = 1472914
$ws.Range("J65").Value = 101558.6
$ws.Range("K65").Value = 7364570
$ws.Range("L65").Value = 507793
$ws.Range("M65").Value = -7361450
$ws.Range("N65").Value = -514033
$ws.Range("H86").Value = 25130318
$ws.Range("I86").Value = 4700
$ws.Range("J86").Value = 50255936
$ws.Range("K86").Value = 4700
$ws.Range("L86").Value = 50255936
$ws.Range("M86").Value = -3577
$ws.Range("N86").Value = -50258182
$ws.Range("H89").Value = 25130318
$ws.Range("I89").Value = 4700
$ws.Range("J89").Value = 50255936
$ws.Range("K89").Value = 23500
$ws.Range("L89").Value = 251279680
$ws.Range("M89").Value = -17884
$ws.Range("N89").Value = -251290912
$ws.Range("H112").Value = 64046.438
$ws.Range("I112").Value = 112364.11
$ws.Range("K112").Value = 337092.33
$ws.Range("M112").Value = -335984.33
$ws.Range("H127").Value = 978.25
$ws.Range("I127").Value = 878.3077
$ws.Range("K127").Value = 2634.9231
$ws.Range("M127").Value = 2325.0769
$ws.Range("H129").Value = 369.75
$ws.Range("I129").Value = 369.75
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1109.25
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3890.75
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 4054.4688
$ws.Range("I132").Value = 3857.1738
$ws.Range("J132").Value = 4558.6665
$ws.Range("K132").Value = 11571.5214
$ws.Range("L132").Value = 13675.9995
$ws.Range("M132").Value = -9041.5214
$ws.Range("N132").Value = -18735.9995
$ws.Range("H135").Value = 29414464
$ws.Range("I135").Value = 37040056
$ws.Range("J135").Value = 1456.1428
$ws.Range("K135").Value = 333360504
$ws.Range("L135").Value = 13105.2852
$ws.Range("M135").Value = -333357969
$ws.Range("N135").Value = -18175.2852
$ws.Range("H137").Value = 1843.5862
$ws.Range("I137").Value = 1421.05
$ws.Range("J137").Value = 2782.5557
$ws.Range("K137").Value = 4263.15
$ws.Range("L137").Value = 8347.667099999999
$ws.Range("M137").Value = -1713.15
$ws.Range("N137").Value = -13447.6671
$ws.Range("H141").Value = 930
$ws.Range("I141").Value = 930
$ws.Range("K141").Value = 2790
$ws.Range("M141").Value = 2390

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2464.9558
$ws.Range("I32").Value = 1400.2931
$ws.Range("K32").Value = 1400.2931
$ws.Range("M32").Value = -1113.2931
$ws.Range("H45").Value = 2501.125
$ws.Range("I45").Value = 2001.5
$ws.Range("K45").Value = 2001.5
$ws.Range("M45").Value = -1624.5
$ws.Range("H61").Value = 111115160
$ws.Range("I61").Value = 166667470
$ws.Range("J61").Value = 10526
$ws.Range("K61").Value = 166667470
$ws.Range("L61").Value = 10526
$ws.Range("M61").Value = -166667258
$ws.Range("N61").Value = -10950
$ws.Range("H74").Value = 1853.5952
$ws.Range("I74").Value = 1637.2778
$ws.Range("J74").Value = 3151.5
$ws.Range("K74").Value = 1637.2778
$ws.Range("L74").Value = 3151.5
$ws.Range("M74").Value = -763.2778000000001
$ws.Range("N74").Value = -4899.5
$ws.Range("H77").Value = 1853.5952
$ws.Range("I77").Value = 1637.2778
$ws.Range("J77").Value = 3151.5
$ws.Range("K77").Value = 8186.389
$ws.Range("L77").Value = 15757.5
$ws.Range("M77").Value = -3818.389
$ws.Range("N77").Value = -24493.5
$ws.Range("H88").Value = 15153854
$ws.Range("I88").Value = 41668268
$ws.Range("J88").Value = 2761.2856
$ws.Range("K88").Value = 41668268
$ws.Range("L88").Value = 2761.2856
$ws.Range("M88").Value = -41667862
$ws.Range("N88").Value = -3573.2856
$ws.Range("H91").Value = 15153854
$ws.Range("I91").Value = 41668268
$ws.Range("J91").Value = 2761.2856
$ws.Range("K91").Value = 41668268
$ws.Range("L91").Value = 2761.2856
$ws.Range("M91").Value = -41666864
$ws.Range("N91").Value = -5569.2856
$ws.Range("H97").Value = 699.6
$ws.Range("I97").Value = 705.4167
$ws.Range("J97").Value = 676.3333
$ws.Range("K97").Value = 705.4167
$ws.Range("L97").Value = 676.3333
$ws.Range("M97").Value = -209.4167
$ws.Range("N97").Value = -1668.3333
$ws.Range("H102").Value = 3377580.5
$ws.Range("I102").Value = 3954037.5
$ws.Range("K102").Value = 3954037.5
$ws.Range("M102").Value = -3952415.5
$ws.Range("H122").Value = 8549859
$ws.Range("I122").Value = 10103653
$ws.Range("J122").Value = 3995.6667
$ws.Range("K122").Value = 30310959
$ws.Range("L122").Value = 11987.0001
$ws.Range("M122").Value = -30308509
$ws.Range("N122").Value = -16887.0001
$ws.Range("H132").Value = 71431120
$ws.Range("I132").Value = 125002420
$ws.Range("K132").Value = 375007260
$ws.Range("M132").Value = -375004730
$ws.Range("H136").Value = 111115160
$ws.Range("I136").Value = 166667470
$ws.Range("J136").Value = 10526
$ws.Range("K136").Value = 500002410
$ws.Range("L136").Value = 31578
$ws.Range("M136").Value = -499999860
$ws.Range("N136").Value = -36678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 515.25
$ws.Range("I20").Value = 564.6667
$ws.Range("J20").Value = 367
$ws.Range("K20").Value = 564.6667
$ws.Range("L20").Value = 367
$ws.Range("M20").Value = -317.6667
$ws.Range("N20").Value = -861
$ws.Range("H38").Value = 51924.7
$ws.Range("J38").Value = 51924.7
$ws.Range("L38").Value = 51924.7
$ws.Range("N38").Value = -52756.7
$ws.Range("H86").Value = 6119.3
$ws.Range("I86").Value = 6587.375
$ws.Range("K86").Value = 6587.375
$ws.Range("M86").Value = -5464.375
$ws.Range("H89").Value = 6119.3
$ws.Range("I89").Value = 6587.375
$ws.Range("K89").Value = 32936.875
$ws.Range("M89").Value = -27320.875
$ws.Range("H94").Value = 5669
$ws.Range("I94").Value = 7094.1
$ws.Range("J94").Value = 2818.8
$ws.Range("K94").Value = 7094.1
$ws.Range("L94").Value = 2818.8
$ws.Range("M94").Value = -6643.1
$ws.Range("N94").Value = -3720.8
$ws.Range("H99").Value = 1133.1428
$ws.Range("I99").Value = 1133.1428
$ws.Range("K99").Value = 1133.1428
$ws.Range("M99").Value = 364.8571999999999
$ws.Range("H107").Value = 55759444
$ws.Range("I107").Value = 219285.72
$ws.Range("J107").Value = 250150000
$ws.Range("K107").Value = 219285.72
$ws.Range("L107").Value = 250150000
$ws.Range("M107").Value = -217365.72
$ws.Range("N107").Value = -250153840
$ws.Range("H134").Value = 3341.75
$ws.Range("I134").Value = 3251.375
$ws.Range("K134").Value = 9754.125
$ws.Range("M134").Value = -7219.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H16").Value = 2032.5
$ws.Range("I16").Value = 2149.25
$ws.Range("K16").Value = 2149.25
$ws.Range("M16").Value = -1862.25
$ws.Range("H44").Value = 10067.5
$ws.Range("I44").Value = 10064
$ws.Range("J44").Value = 10071
$ws.Range("K44").Value = 10064
$ws.Range("L44").Value = 10071
$ws.Range("M44").Value = -9622
$ws.Range("N44").Value = -10955
$ws.Range("H47").Value = 42499.75
$ws.Range("I47").Value = 39999.668
$ws.Range("K47").Value = 39999.668
$ws.Range("M47").Value = -39433.668
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H58").Value = 1985.5883
$ws.Range("I58").Value = 1980.2727
$ws.Range("J58").Value = 1995.3334
$ws.Range("K58").Value = 1980.2727
$ws.Range("L58").Value = 1995.3334
$ws.Range("M58").Value = -1777.2727
$ws.Range("N58").Value = -2401.3334
$ws.Range("H99").Value = 2614.611
$ws.Range("I99").Value = 2191.5454
$ws.Range("J99").Value = 3279.4285
$ws.Range("K99").Value = 2191.5454
$ws.Range("L99").Value = 3279.4285
$ws.Range("M99").Value = -693.5454
$ws.Range("N99").Value = -6275.4285
$ws.Range("H113").Value = 2032.5
$ws.Range("I113").Value = 2149.25
$ws.Range("K113").Value = 2149.25
$ws.Range("M113").Value = 20.75
$ws.Range("H126").Value = 2614.611
$ws.Range("I126").Value = 2191.5454
$ws.Range("J126").Value = 3279.4285
$ws.Range("K126").Value = 6574.6362
$ws.Range("L126").Value = 9838.2855
$ws.Range("M126").Value = -4104.6362
$ws.Range("N126").Value = -14778.2855
$ws.Range("H132").Value = 4298
$ws.Range("I132").Value = 4357.8
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 13073.4
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -10543.4
$ws.Range("N132").Value = -17057
$ws.Range("H136").Value = 1985.5883
$ws.Range("I136").Value = 1980.2727
$ws.Range("J136").Value = 1995.3334
$ws.Range("K136").Value = 5940.8181
$ws.Range("L136").Value = 5986.0002
$ws.Range("M136").Value = -3390.8181
$ws.Range("N136").Value = -11086.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1869.5333
$ws.Range("I8").Value = 1869.5333
$ws.Range("K8").Value = 5608.5999
$ws.Range("M8").Value = -5469.5999
$ws.Range("H12").Value = 336.8889
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 336.8889
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1010.6667
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1356.6667
$ws.Range("H92").Value = 239.2
$ws.Range("J92").Value = 303
$ws.Range("L92").Value = 909
$ws.Range("N92").Value = -3405
$ws.Range("H113").Value = 688.30304
$ws.Range("I113").Value = 443.45456
$ws.Range("J113").Value = 810.7273
$ws.Range("K113").Value = 1330.36368
$ws.Range("L113").Value = 2432.1819
$ws.Range("M113").Value = 839.6363200000001
$ws.Range("N113").Value = -6772.1819
$ws.Range("H132").Value = 1361.75
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1361.75
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12255.75
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -17315.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2754
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H11").Value = 4287250
$ws.Range("J11").Value = 5252062.5
$ws.Range("L11").Value = 5252062.5
$ws.Range("N11").Value = -5252340.5
$ws.Range("H70").Value = 8050
$ws.Range("I70").Value = 7656.353
$ws.Range("K70").Value = 7656.353
$ws.Range("M70").Value = -7386.353
$ws.Range("H73").Value = 8050
$ws.Range("I73").Value = 7656.353
$ws.Range("K73").Value = 7656.353
$ws.Range("M73").Value = -6720.353
$ws.Range("H80").Value = 3231.7
$ws.Range("I80").Value = 3345.6667
$ws.Range("J80").Value = 3060.75
$ws.Range("K80").Value = 3345.6667
$ws.Range("L80").Value = 3060.75
$ws.Range("M80").Value = -2347.6667
$ws.Range("N80").Value = -5056.75
$ws.Range("H83").Value = 3231.7
$ws.Range("I83").Value = 3345.6667
$ws.Range("J83").Value = 3060.75
$ws.Range("K83").Value = 16728.3335
$ws.Range("L83").Value = 15303.75
$ws.Range("M83").Value = -11736.3335
$ws.Range("N83").Value = -25287.75
$ws.Range("H97").Value = 999.05554
$ws.Range("I97").Value = 767.0769
$ws.Range("K97").Value = 767.0769
$ws.Range("M97").Value = -271.0769
$ws.Range("H126").Value = 27503
$ws.Range("I126").Value = 35337.332
$ws.Range("K126").Value = 106011.996
$ws.Range("M126").Value = -103541.996
$ws.Range("H132").Value = 4362.2163
$ws.Range("I132").Value = 3833.5
$ws.Range("J132").Value = 5611.909
$ws.Range("K132").Value = 11500.5
$ws.Range("L132").Value = 16835.727
$ws.Range("M132").Value = -8970.5
$ws.Range("N132").Value = -21895.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45456548
$ws.Range("I7").Value = 55557340
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 55557340
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -55557228
$ws.Range("N7").Value = -3224
$ws.Range("H10").Value = 55699.6
$ws.Range("J10").Value = 55699.6
$ws.Range("L10").Value = 55699.6
$ws.Range("N10").Value = -55979.6
$ws.Range("H12").Value = 1001999.4
$ws.Range("I12").Value = 2500501.5
$ws.Range("J12").Value = 2998
$ws.Range("K12").Value = 2500501.5
$ws.Range("L12").Value = 2998
$ws.Range("M12").Value = -2500331.5
$ws.Range("N12").Value = -3338
$ws.Range("H22").Value = 1749.75
$ws.Range("J22").Value = 1333
$ws.Range("L22").Value = 1333
$ws.Range("N22").Value = -1923
$ws.Range("H27").Value = 1749.75
$ws.Range("J27").Value = 1333
$ws.Range("L27").Value = 1333
$ws.Range("N27").Value = -1547
$ws.Range("H55").Value = 553.5
$ws.Range("J55").Value = 594.8889
$ws.Range("L55").Value = 594.8889
$ws.Range("N55").Value = -940.8889
$ws.Range("H93").Value = 1308.1613
$ws.Range("I93").Value = 1257.6296
$ws.Range("K93").Value = 1257.6296
$ws.Range("M93").Value = -9.629599999999982
$ws.Range("H100").Value = 3624.75
$ws.Range("I100").Value = 3499.6667
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3499.6667
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2958.6667
$ws.Range("N100").Value = -5082
$ws.Range("H122").Value = 5227.1816
$ws.Range("I122").Value = 5208.1665
$ws.Range("K122").Value = 15624.4995
$ws.Range("M122").Value = -13174.4995
$ws.Range("H126").Value = 45456548
$ws.Range("I126").Value = 55557340
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 166672020
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -166669550
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 6669.905
$ws.Range("I132").Value = 3732.6667
$ws.Range("J132").Value = 8872.833000000001
$ws.Range("K132").Value = 11198.0001
$ws.Range("L132").Value = 26618.499
$ws.Range("M132").Value = -8668.000100000001
$ws.Range("N132").Value = -31678.499
$ws.Range("H135").Value = 42365.4
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 42365.4
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 42365.4
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -52505.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 131011.22
$ws.Range("I4").Value = 3683.6667
$ws.Range("J4").Value = 385666.34
$ws.Range("K4").Value = 3683.6667
$ws.Range("L4").Value = 385666.34
$ws.Range("M4").Value = -3570.6667
$ws.Range("N4").Value = -385892.34
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H13").Value = 2500
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2500
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2500
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -2780
$ws.Range("H96").Value = 1965.6666
$ws.Range("I96").Value = 1965.6666
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1965.6666
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -592.6666
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 2173.5
$ws.Range("I122").Value = 1898
$ws.Range("K122").Value = 5694
$ws.Range("M122").Value = -3244
$ws.Range("H126").Value = 1878.1765
$ws.Range("I126").Value = 1746.1
$ws.Range("J126").Value = 2066.8572
$ws.Range("K126").Value = 5238.299999999999
$ws.Range("L126").Value = 6200.571599999999
$ws.Range("M126").Value = -2768.299999999999
$ws.Range("N126").Value = -11140.5716
$ws.Range("H132").Value = 3784.0605
$ws.Range("I132").Value = 4403.1577
$ws.Range("J132").Value = 2943.8572
$ws.Range("K132").Value = 13209.4731
$ws.Range("L132").Value = 8831.571599999999
$ws.Range("M132").Value = -10679.4731
$ws.Range("N132").Value = -13891.5716
$ws.Range("H136").Value = 4391.6
$ws.Range("I136").Value = 2528.7778
$ws.Range("J136").Value = 5915.727
$ws.Range("K136").Value = 7586.3334
$ws.Range("L136").Value = 17747.181
$ws.Range("M136").Value = -5036.3334
$ws.Range("N136").Value = -22847.181
